$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-29 Thursday", "2026-01-30 Friday"),
    @("83×59=", "98×65="),
    @("81×76=", "54×25="),
    @("91×60=", "38×86="),
    @("31×34=", "83×22="),
    @("96×15=", "42×16="),
    @("39×80=", "48×59="),
    @("49×33=", "36×36="),
    @("36×79=", "11×11="),
    @("37×48=", "66×21="),
    @("17×14=", "40×71="),
    @("35×98=", "77×27="),
    @("27×13=", "20×33="),
    @("61×57=", "28×16="),
    @("83×92=", "63×14="),
    @("86×57=", "83×88="),
    @("76×82=", "86×30="),
    @("42×68=", "40×12="),
    @("29×27=", "49×99="),
    @("77×70=", "83×46="),
    @("62×36=", "80×72="),
    @("32×17=", "15×51="),
    @("82×35=", "45×96="),
    @("75×36=", "80×91="),
    @("58×47=", "77×81="),
    @("82×41=", "98×96=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
